$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 74.30768999999999
$ws.Range("I2").Value = 51.454544
$ws.Range("K2").Value = 51.454544
$ws.Range("M2").Value = 61.545456
$ws.Range("H9").Value = 229.82608
$ws.Range("I9").Value = 107.666664
$ws.Range("J9").Value = 669.6
$ws.Range("K9").Value = 107.666664
$ws.Range("L9").Value = 669.6
$ws.Range("M9").Value = 61.333336
$ws.Range("N9").Value = -1007.6
$ws.Range("H40").Value = 2333.1482
$ws.Range("I40").Value = 2437.2917
$ws.Range("K40").Value = 2437.2917
$ws.Range("M40").Value = -2262.2917
$ws.Range("H41").Value = 507.66666
$ws.Range("I41").Value = 321.25
$ws.Range("K41").Value = 321.25
$ws.Range("M41").Value = 118.75
$ws.Range("H64").Value = 3749.55
$ws.Range("I64").Value = 2999.3076
$ws.Range("K64").Value = 2999.3076
$ws.Range("M64").Value = -2751.3076
$ws.Range("H67").Value = 3749.55
$ws.Range("I67").Value = 2999.3076
$ws.Range("K67").Value = 2999.3076
$ws.Range("M67").Value = -2141.3076
$ws.Range("H69").Value = 7224.4443
$ws.Range("J69").Value = 7224.4443
$ws.Range("L69").Value = 21673.3329
$ws.Range("N69").Value = -23421.3329
$ws.Range("H70").Value = 3374.75
$ws.Range("I70").Value = 2785.4285
$ws.Range("J70").Value = 4199.8
$ws.Range("K70").Value = 8356.2855
$ws.Range("L70").Value = 12599.4
$ws.Range("M70").Value = -8086.2855
$ws.Range("N70").Value = -13139.4
$ws.Range("H72").Value = 7224.4443
$ws.Range("J72").Value = 7224.4443
$ws.Range("L72").Value = 65019.9987
$ws.Range("N72").Value = -73755.9987
$ws.Range("H73").Value = 3374.75
$ws.Range("I73").Value = 2785.4285
$ws.Range("J73").Value = 4199.8
$ws.Range("K73").Value = 8356.2855
$ws.Range("L73").Value = 12599.4
$ws.Range("M73").Value = -7420.2855
$ws.Range("N73").Value = -14471.4
$ws.Range("H80").Value = 2239.1428
$ws.Range("I80").Value = 1231.375
$ws.Range("K80").Value = 3694.125
$ws.Range("M80").Value = -2696.125
$ws.Range("H83").Value = 2239.1428
$ws.Range("I83").Value = 1231.375
$ws.Range("K83").Value = 11082.375
$ws.Range("M83").Value = -6090.375
$ws.Range("H87").Value = 52000
$ws.Range("J87").Value = 66666.664
$ws.Range("L87").Value = 66666.664
$ws.Range("N87").Value = -69162.664
$ws.Range("H90").Value = 52000
$ws.Range("J90").Value = 66666.664
$ws.Range("L90").Value = 199999.992
$ws.Range("N90").Value = -212479.992
$ws.Range("H92").Value = 1398
$ws.Range("I92").Value = 1666.6666
$ws.Range("K92").Value = 1666.6666
$ws.Range("M92").Value = -418.6666
$ws.Range("H94").Value = 1997.4286
$ws.Range("I94").Value = 1497
$ws.Range("K94").Value = 1497
$ws.Range("M94").Value = -1046
$ws.Range("H99").Value = 463.18182
$ws.Range("I99").Value = 299.625
$ws.Range("J99").Value = 899.3333
$ws.Range("K99").Value = 898.875
$ws.Range("L99").Value = 2697.9999
$ws.Range("M99").Value = 599.125
$ws.Range("N99").Value = -5693.9999
$ws.Range("H100").Value = 2293.5881
$ws.Range("I100").Value = 1320
$ws.Range("J100").Value = 2699.25
$ws.Range("K100").Value = 1320
$ws.Range("L100").Value = 2699.25
$ws.Range("M100").Value = -779
$ws.Range("N100").Value = -3781.25
$ws.Range("H101").Value = 1713.7142
$ws.Range("I101").Value = 383.33334
$ws.Range("J101").Value = 9696
$ws.Range("K101").Value = 1150.00002
$ws.Range("L101").Value = 29088
$ws.Range("M101").Value = 471.9999800000001
$ws.Range("N101").Value = -32332
$ws.Range("H116").Value = 17744.088
$ws.Range("I116").Value = 6729.4165
$ws.Range("J116").Value = 29760.092
$ws.Range("K116").Value = 6729.4165
$ws.Range("L116").Value = 29760.092
$ws.Range("M116").Value = -3287.4165
$ws.Range("N116").Value = -36644.092
$ws.Range("H118").Value = 580
$ws.Range("I118").Value = 580
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1740
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -83
$ws.Range("N118").ClearContents()
$ws.Range("H137").Value = 45464428
$ws.Range("I137").Value = 125003384
$ws.Range("J137").Value = 13598.5
$ws.Range("K137").Value = 375010152
$ws.Range("L137").Value = 40795.5
$ws.Range("M137").Value = -375007602
$ws.Range("N137").Value = -45895.5
$ws.Range("H138").Value = 5923
$ws.Range("I138").Value = 7884.3335
$ws.Range("J138").Value = 5649.3257
$ws.Range("K138").Value = 23653.0005
$ws.Range("L138").Value = 16947.9771
$ws.Range("M138").Value = -18513.0005
$ws.Range("N138").Value = -27227.9771
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2502.889
$ws.Range("I2").Value = 2844.3333
$ws.Range("K2").Value = 2844.3333
$ws.Range("M2").Value = -2731.3333
$ws.Range("H5").Value = 2848.5557
$ws.Range("I5").Value = 2067.4
$ws.Range("J5").Value = 3825
$ws.Range("K5").Value = 2067.4
$ws.Range("L5").Value = 3825
$ws.Range("M5").Value = -1955.4
$ws.Range("N5").Value = -4049
$ws.Range("H32").Value = 154041.66
$ws.Range("I32").Value = 181454.86
$ws.Range("J32").Value = 33845.31
$ws.Range("K32").Value = 181454.86
$ws.Range("L32").Value = 33845.31
$ws.Range("M32").Value = -181167.86
$ws.Range("N32").Value = -34419.31
$ws.Range("H61").Value = 13404769
$ws.Range("I61").Value = 5961.1665
$ws.Range("K61").Value = 5961.1665
$ws.Range("M61").Value = -5749.1665
$ws.Range("H74").Value = 1301741.6
$ws.Range("I74").Value = 1640946.9
$ws.Range("K74").Value = 1640946.9
$ws.Range("M74").Value = -1640072.9
$ws.Range("H77").Value = 1301741.6
$ws.Range("I77").Value = 1640946.9
$ws.Range("K77").Value = 8204734.5
$ws.Range("M77").Value = -8200366.5
$ws.Range("H97").Value = 29413234
$ws.Range("I97").Value = 1285.7097
$ws.Range("J97").Value = 333336700
$ws.Range("K97").Value = 1285.7097
$ws.Range("L97").Value = 333336700
$ws.Range("M97").Value = -789.7097000000001
$ws.Range("N97").Value = -333337692
$ws.Range("H110").Value = 1372.2
$ws.Range("I110").Value = 964.9231
$ws.Range("J110").Value = 2128.5715
$ws.Range("K110").Value = 964.9231
$ws.Range("L110").Value = 2128.5715
$ws.Range("M110").Value = 1080.0769
$ws.Range("N110").Value = -6218.5715
$ws.Range("H116").Value = 2502.889
$ws.Range("I116").Value = 2844.3333
$ws.Range("K116").Value = 2844.3333
$ws.Range("M116").Value = -550.3332999999998
$ws.Range("H122").Value = 62501600
$ws.Range("I122").Value = 83334300
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 250002900
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -250000450
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 2086433.4
$ws.Range("J132").Value = 3928.4285
$ws.Range("L132").Value = 11785.2855
$ws.Range("N132").Value = -16845.2855
$ws.Range("H136").Value = 13404769
$ws.Range("I136").Value = 5961.1665
$ws.Range("K136").Value = 17883.4995
$ws.Range("M136").Value = -15333.4995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2502.889
$ws.Range("I3").Value = 2844.3333
$ws.Range("K3").Value = 2844.3333
$ws.Range("M3").Value = -2730.3333
$ws.Range("H4").Value = 2848.5557
$ws.Range("I4").Value = 2067.4
$ws.Range("J4").Value = 3825
$ws.Range("K4").Value = 2067.4
$ws.Range("L4").Value = 3825
$ws.Range("M4").Value = -1952.4
$ws.Range("N4").Value = -4055
$ws.Range("H20").Value = 48052.332
$ws.Range("I20").Value = 79996.14
$ws.Range("K20").Value = 79996.14
$ws.Range("M20").Value = -79749.14
$ws.Range("H22").Value = 173.44444
$ws.Range("I22").Value = 144.5
$ws.Range("J22").Value = 196.6
$ws.Range("K22").Value = 144.5
$ws.Range("L22").Value = 196.6
$ws.Range("M22").Value = 28.5
$ws.Range("N22").Value = -542.6
$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9761
$ws.Range("H99").Value = 6238.2905
$ws.Range("I99").Value = 8775.9
$ws.Range("K99").Value = 8775.9
$ws.Range("M99").Value = -7277.9
$ws.Range("H107").Value = 22728482
$ws.Range("J107").Value = 1416.5
$ws.Range("L107").Value = 1416.5
$ws.Range("N107").Value = -5256.5
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 5218343.5
$ws.Range("I134").Value = 6196.278
$ws.Range("K134").Value = 18588.834
$ws.Range("M134").Value = -16053.834
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 73020.5
$ws.Range("J22").Value = 127561.75
$ws.Range("L22").Value = 127561.75
$ws.Range("N22").Value = -128261.75
$ws.Range("H31").Value = 2648559.5
$ws.Range("I31").Value = 3089322.2
$ws.Range("J31").Value = 3983.3333
$ws.Range("K31").Value = 3089322.2
$ws.Range("L31").Value = 3983.3333
$ws.Range("M31").Value = -3089027.2
$ws.Range("N31").Value = -4573.3333
$ws.Range("H34").Value = 2648559.5
$ws.Range("I34").Value = 3089322.2
$ws.Range("J34").Value = 3983.3333
$ws.Range("K34").Value = 3089322.2
$ws.Range("L34").Value = 3983.3333
$ws.Range("M34").Value = -3089120.2
$ws.Range("N34").Value = -4387.3333
$ws.Range("H75").Value = 19999
$ws.Range("J75").Value = 19999
$ws.Range("L75").Value = 19999
$ws.Range("N75").Value = -21995
$ws.Range("H78").Value = 19999
$ws.Range("J78").Value = 19999
$ws.Range("L78").Value = 59997
$ws.Range("N78").Value = -69981
$ws.Range("H99").Value = 38632.523
$ws.Range("I99").Value = 37149.875
$ws.Range("J99").Value = 42021.43
$ws.Range("K99").Value = 37149.875
$ws.Range("L99").Value = 42021.43
$ws.Range("M99").Value = -35651.875
$ws.Range("N99").Value = -45017.43
$ws.Range("H105").Value = 5591.375
$ws.Range("I105").Value = 5905.316
$ws.Range("K105").Value = 5905.316
$ws.Range("M105").Value = -4158.316
$ws.Range("H122").Value = 28951.875
$ws.Range("I122").Value = 9723.4
$ws.Range("J122").Value = 60999.332
$ws.Range("K122").Value = 29170.2
$ws.Range("L122").Value = 182997.996
$ws.Range("M122").Value = -26720.2
$ws.Range("N122").Value = -187897.996
$ws.Range("H126").Value = 38632.523
$ws.Range("I126").Value = 37149.875
$ws.Range("J126").Value = 42021.43
$ws.Range("K126").Value = 111449.625
$ws.Range("L126").Value = 126064.29
$ws.Range("M126").Value = -108979.625
$ws.Range("N126").Value = -131004.29
$ws.Range("H132").Value = 9473.762000000001
$ws.Range("I132").Value = 5400.0713
$ws.Range("J132").Value = 17621.143
$ws.Range("K132").Value = 16200.2139
$ws.Range("L132").Value = 52863.429
$ws.Range("M132").Value = -13670.2139
$ws.Range("N132").Value = -57923.429
$ws.Range("H134").Value = 3149.8462
$ws.Range("I134").Value = 2680.6316
$ws.Range("J134").Value = 4423.4287
$ws.Range("K134").Value = 8041.8948
$ws.Range("L134").Value = 13270.2861
$ws.Range("M134").Value = -5506.8948
$ws.Range("N134").Value = -18340.2861
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2107826.8
$ws.Range("I5").Value = 1701317.1
$ws.Range("K5").Value = 5103951.300000001
$ws.Range("M5").Value = -5103839.300000001
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 15000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -15346
$ws.Range("H113").Value = 2527.875
$ws.Range("J113").Value = 3318
$ws.Range("L113").Value = 9954
$ws.Range("N113").Value = -14294
$ws.Range("H121").Value = 3127416.5
$ws.Range("J121").Value = 6254390.5
$ws.Range("L121").Value = 18763171.5
$ws.Range("N121").Value = -18765791.5
$ws.Range("H122").Value = 1614394.2
$ws.Range("J122").Value = 1801
$ws.Range("L122").Value = 16209
$ws.Range("N122").Value = -21109
$ws.Range("H127").Value = 9690.454
$ws.Range("J127").Value = 9690.454
$ws.Range("L127").Value = 29071.362
$ws.Range("N127").Value = -38991.362
$ws.Range("H129").Value = 6312.0625
$ws.Range("I129").Value = 6008.75
$ws.Range("J129").Value = 6413.1665
$ws.Range("K129").Value = 18026.25
$ws.Range("L129").Value = 19239.4995
$ws.Range("M129").Value = -13026.25
$ws.Range("N129").Value = -29239.4995
$ws.Range("H131").Value = 8244.076999999999
$ws.Range("I131").Value = 1687
$ws.Range("K131").Value = 5061
$ws.Range("M131").Value = -21
$ws.Range("H135").Value = 2107826.8
$ws.Range("I135").Value = 1701317.1
$ws.Range("K135").Value = 15311853.9
$ws.Range("M135").Value = -15309318.9
$ws.Range("H140").Value = 3807.3667
$ws.Range("I140").Value = 2491.3809
$ws.Range("K140").Value = 7474.1427
$ws.Range("M140").Value = -2294.1427
$ws.Range("H141").Value = 7580.353
$ws.Range("I141").Value = 2169.6365
$ws.Range("K141").Value = 6508.9095
$ws.Range("M141").Value = -1328.9095
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5941443.5
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326
$ws.Range("H102").Value = 957.6875
$ws.Range("I102").Value = 743.6
$ws.Range("J102").Value = 4169
$ws.Range("K102").Value = 743.6
$ws.Range("L102").Value = 4169
$ws.Range("M102").Value = 878.4
$ws.Range("N102").Value = -7413
$ws.Range("H104").Value = 66435.5
$ws.Range("J104").Value = 66435.5
$ws.Range("L104").Value = 66435.5
$ws.Range("N104").Value = -73423.5
$ws.Range("H113").Value = 4263.5
$ws.Range("I113").Value = 4263.5
$ws.Range("K113").Value = 4263.5
$ws.Range("M113").Value = -2093.5
$ws.Range("H122").Value = 56645.1
$ws.Range("I122").Value = 96412.55
$ws.Range("K122").Value = 289237.65
$ws.Range("M122").Value = -286787.65
$ws.Range("H126").Value = 12945.6875
$ws.Range("I126").Value = 13937.929
$ws.Range("K126").Value = 41813.787
$ws.Range("M126").Value = -39343.787
$ws.Range("H132").Value = 29087.25
$ws.Range("I132").Value = 27502.75
$ws.Range("K132").Value = 82508.25
$ws.Range("M132").Value = -79978.25
$ws.Range("H134").Value = 55220.637
$ws.Range("J134").Value = 55220.637
$ws.Range("L134").Value = 165661.911
$ws.Range("N134").Value = -170731.911
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H22").Value = 2642.818
$ws.Range("I22").Value = 1691
$ws.Range("J22").Value = 2999.75
$ws.Range("K22").Value = 1691
$ws.Range("L22").Value = 2999.75
$ws.Range("M22").Value = -1396
$ws.Range("N22").Value = -3589.75
$ws.Range("H27").Value = 2642.818
$ws.Range("I27").Value = 1691
$ws.Range("J27").Value = 2999.75
$ws.Range("K27").Value = 1691
$ws.Range("L27").Value = 2999.75
$ws.Range("M27").Value = -1584
$ws.Range("N27").Value = -3213.75
$ws.Range("H40").Value = 4988.9165
$ws.Range("I40").Value = 4533.364
$ws.Range("K40").Value = 4533.364
$ws.Range("M40").Value = -4397.364
$ws.Range("H46").Value = 5337.2
$ws.Range("I46").Value = 2310.3333
$ws.Range("K46").Value = 2310.3333
$ws.Range("M46").Value = -2122.3333
$ws.Range("H47").Value = 4000
$ws.Range("I47").Value = 4000
$ws.Range("K47").Value = 4000
$ws.Range("M47").Value = -3510
$ws.Range("H52").Value = 4000
$ws.Range("I52").Value = 4000
$ws.Range("K52").Value = 4000
$ws.Range("M52").Value = -3767
$ws.Range("H74").Value = 41994
$ws.Range("I74").Value = 41994
$ws.Range("K74").Value = 41994
$ws.Range("M74").Value = -40996
$ws.Range("H77").Value = 41994
$ws.Range("I77").Value = 41994
$ws.Range("K77").Value = 125982
$ws.Range("M77").Value = -120990
$ws.Range("H93").Value = 4537
$ws.Range("I93").Value = 3854.25
$ws.Range("J93").Value = 9999
$ws.Range("K93").Value = 3854.25
$ws.Range("L93").Value = 9999
$ws.Range("M93").Value = -2606.25
$ws.Range("N93").Value = -12495
$ws.Range("H122").Value = 7486
$ws.Range("I122").Value = 8076.5
$ws.Range("K122").Value = 24229.5
$ws.Range("M122").Value = -21779.5
$ws.Range("H132").Value = 19478426
$ws.Range("I132").Value = 23372710
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 70118130
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -70115600
$ws.Range("N132").Value = -26060
$ws.Range("H134").Value = 60420
$ws.Range("J134").Value = 60420
$ws.Range("L134").Value = 60420
$ws.Range("N134").Value = -70560
$ws.Range("H136").Value = 9270002
$ws.Range("I136").Value = 6951883
$ws.Range("J136").Value = 13906239
$ws.Range("K136").Value = 20855649
$ws.Range("L136").Value = 41718717
$ws.Range("M136").Value = -20853099
$ws.Range("N136").Value = -41723817
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 23068.625
$ws.Range("I51").Value = 19829.5
$ws.Range("J51").Value = 26307.75
$ws.Range("K51").Value = 19829.5
$ws.Range("L51").Value = 26307.75
$ws.Range("M51").Value = -19319.5
$ws.Range("N51").Value = -27327.75
$ws.Range("H52").Value = 18896
$ws.Range("J52").Value = 23437.6
$ws.Range("L52").Value = 23437.6
$ws.Range("N52").Value = -23889.6
$ws.Range("H62").Value = 70750.75
$ws.Range("J62").Value = 84334.336
$ws.Range("L62").Value = 84334.336
$ws.Range("N62").Value = -85582.336
$ws.Range("H65").Value = 70750.75
$ws.Range("J65").Value = 84334.336
$ws.Range("L65").Value = 421671.68
$ws.Range("N65").Value = -427911.68
$ws.Range("H81").Value = 4068.4
$ws.Range("I81").Value = 2924.625
$ws.Range("K81").Value = 5849.25
$ws.Range("M81").Value = -4788.25
$ws.Range("H84").Value = 4068.4
$ws.Range("I84").Value = 2924.625
$ws.Range("K84").Value = 29246.25
$ws.Range("M84").Value = -23942.25
$ws.Range("H122").Value = 70762.94
$ws.Range("I122").Value = 5076.933
$ws.Range("J122").Value = 563408
$ws.Range("K122").Value = 15230.799
$ws.Range("L122").Value = 1690224
$ws.Range("M122").Value = -12780.799
$ws.Range("N122").Value = -1695124
$ws.Range("H132").Value = 55558224
$ws.Range("I132").Value = 55558224
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 166674672
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -166672142
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 14807834
$ws.Range("I136").Value = 2901272.5
$ws.Range("J136").Value = 50527516
$ws.Range("K136").Value = 8703817.5
$ws.Range("L136").Value = 151582548
$ws.Range("M136").Value = -8701267.5
$ws.Range("N136").Value = -151587648
